$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to be treated as plain text so Excel
# doesn't silently coerce values that look numeric (e.g. '240.95') into
# floating point numbers - the source data is text in every row.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "36.549.93"
$ws.Range("E2").Value = "  -1.43%  "

$ws.Range("D3").Value = "2.058.23"
$ws.Range("E3").Value = "  +0.65%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "240.95"
$ws.Range("E5").Value = "  -2.69%  "

$ws.Range("E6").Value = "  +1.68%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "52.11"
$ws.Range("E8").Value = "  -7.70%  "

$ws.Range("D9").Value = "58.72"
$ws.Range("E9").Value = "  -2.13%  "

$ws.Range("D10").Value = "0.359"
$ws.Range("E10").Value = "  -6.49%  "

$ws.Range("D11").Value = "0.0749"
$ws.Range("E11").Value = "  -3.96%  "

$ws.Range("E12").Value = "  -1.24%  "

$ws.Range("D13").Value = "0.883"
$ws.Range("E13").Value = "  -1.74%  "

$ws.Range("D14").Value = "14.42"
$ws.Range("E14").Value = "  -10.03%  "

$ws.Range("D15").Value = "2.364.43"
$ws.Range("E15").Value = "  +0.96%  "

$ws.Range("D16").Value = "5.38"
$ws.Range("E16").Value = "  -5.87%  "

$ws.Range("D17").Value = "2.054.66"
$ws.Range("E17").Value = "  +0.51%  "

$ws.Range("D18").Value = "36.473.22"
$ws.Range("E18").Value = "  -1.74%  "

$ws.Range("D19").Value = "16.24"
$ws.Range("E19").Value = "  -13.94%  "

$ws.Range("D20").Value = "71.45"
$ws.Range("E20").Value = "  -4.59%  "

$ws.Range("D21").Value = "0.0₃0858"
$ws.Range("E21").Value = "  -3.74%  "

$ws.Range("D22").Value = "5.24"
$ws.Range("E22").Value = "  -3.26%  "

$ws.Range("D23").Value = "235.17"
$ws.Range("E23").Value = "  -0.80%  "

$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("D25").Value = "2.35"
$ws.Range("E25").Value = "  -5.27%  "

$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "2.12"
$ws.Range("E26").Value = "  -2.50%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "9.22"
$ws.Range("E27").Value = "  -3.43%  "

$ws.Range("D28").Value = "162.93"
$ws.Range("E28").Value = "  -4.81%  "

$ws.Range("D29").Value = "20.21"
$ws.Range("E29").Value = "  +0.52%  "

$ws.Range("E30").Value = "  -1.21%  "

$ws.Range("D31").Value = "5.02"
$ws.Range("E31").Value = "  -2.37%  "

$ws.Range("E32").Value = "  -4.17%  "

$ws.Range("D33").Value = "4.53"
$ws.Range("E33").Value = "  -3.09%  "

$ws.Range("D34").Value = "0.0593"
$ws.Range("E34").Value = "  -4.83%  "

$ws.Range("E35").Value = "  +0.11%  "

$ws.Range("D36").Value = "1.83"
$ws.Range("E36").Value = "  -2.15%  "

$ws.Range("D37").Value = "2.24"
$ws.Range("E37").Value = "  -1.22%  "

$ws.Range("E38").Value = "  -7.63%  "

$ws.Range("E39").Value = "  -6.84%  "

$ws.Range("D40").Value = "4.85"
$ws.Range("E40").Value = "  -5.86%  "

$ws.Range("E41").Value = "  -6.18%  "

$ws.Range("D42").Value = "0.0214"
$ws.Range("E42").Value = "  -4.07%  "

$ws.Range("E43").Value = "  -3.00%  "

$ws.Range("D44").Value = "0.0927"
$ws.Range("E44").Value = "  -6.41%  "

$ws.Range("D45").Value = "93.27"
$ws.Range("E45").Value = "  -6.46%  "

$ws.Range("D46").Value = "1.383.75"
$ws.Range("E46").Value = "  +7.43%  "

$ws.Range("D47").Value = "15.49"
$ws.Range("E47").Value = "  -10.43%  "

$ws.Range("D48").Value = "7.24"
$ws.Range("E48").Value = "  +6.13%  "

$ws.Range("E49").Value = "  -4.13%  "

$ws.Range("E50").Value = "  -0.63%  "

$ws.Range("D51").Value = "2.247.24"
$ws.Range("E51").Value = "  +0.86%  "

# Restore the original (unstyled / General) appearance now that the
# text values are safely stored.
$ws.Range("D2:E51").Style = "Normal"
